$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "202 946 8050"
$ws.Range("A3").Value = "202 712 5974"
$ws.Range("A4").Value = "202 205 8823"
